# Update scripts wuth new tpm
# Writes the recomputed TPM-derived values into columns M:T for rows 2-7
# of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ M = 240.859306;          N = 481.718612;          O = 0.770495783518506;  P = 0.7103085338816119; Q = 6.739323668315333;  R = 40.435942009892;    S = 0.770495783518506;  T = 0.7103085338816119 }
    3 = @{                          O = 0.133770027207319;  P = 0.1849808797181794; Q = 1.170051192691334;                            S = 0.133770027207319;  T = 0.1849808797181794 }
    4 = @{ M = 3.307267666666667;   N = 9.921803000000001;  O = 0.01057976888853842; P = 0.01462999594956937; Q = 0.09253845173588889; R = 0.8328460656230001; S = 0.01057976888853842; T = 0.01462999594956937 }
    5 = @{ M = 18.767532;           N = 37.535064;           O = 0.06003631129389966; P = 0.05534657705729765; Q = 0.525121801204;      R = 3.150730807224;     S = 0.06003631129389966; T = 0.05534657705729765 }
    6 = @{ M = 1.967337333333333;   N = 5.902012;             O = 0.006293404831498911; P = 0.008702693618721296; Q = 0.05504675436577777; R = 0.495420789292; S = 0.006293404831498911; T = 0.008702693618721296 }
    7 = @{ M = 5.884659333333333;   N = 17.653978;            O = 0.01882470426023795; P = 0.02603131977462027; Q = 0.1646547296997778;  R = 1.481892567298;    S = 0.01882470426023795; T = 0.02603131977462027 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
